$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "1.002", "5.910").
# Excel would normally auto-convert such strings to real numbers when assigned
# via .Value, but the source data must stay text (as in the original workbook).
# Temporarily mark the D column as Text format while writing the new values,
# then restore the original ("Normal") cell style so no visual/style diff remains.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.894.20'
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("D3").Value = '1.888.08'
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '331.95'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4624'
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("D8").Value = '0.4112'
$ws.Range("E8").Value = '  +3.63%  '
$ws.Range("D9").Value = '47.54'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '0.07984'
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '0.9941'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").Value = '21.72'
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D13").Value = '1.885.48'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '5.916'
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").Value = '7.067'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("D16").Value = '89.30'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").Value = '17.49'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '28.941.05'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").Value = '5.383'
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").Value = '11.26'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").Value = '2.217'
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").Value = '2.119.88'
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("D27").Value = '157.56'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '19.66'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("D29").Value = '2.119'
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").Value = '5.422'
$ws.Range("D31").Value = '117.94'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").Value = '0.9786'
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("D33").Value = '0.09350'
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("D34").Value = '1.415'
$ws.Range("E34").Value = '  +2.74%  '
$ws.Range("D35").Value = '3.605'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").Value = '5.280'
$ws.Range("D37").Value = '0.06062'
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("D38").Value = '0.02232'
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("D39").Value = '8.294'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = '1.176'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").Value = '0.5775'
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("D44").Value = '0.1818'
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("D45").Value = '1.263'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '2.284'
$ws.Range("E46").Value = '  +11.14%  '
$ws.Range("D47").Value = '0.5484'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").Value = '11.99'
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").Value = '1.908'
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").Value = '0.07010'
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("D51").Value = '110.93'
$ws.Range("E51").Value = '  -0.83%  '

$priceRange.Style = "Normal"
